# Weekly update: insert 3 new rows of "Choclo" price data at row 616,
# pushing the existing rows (616:665) down to (619:668).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 616 (shifts 616:665 -> 619:668).
$ws.Range("616:618").Insert()

# New row 616: Lluteño / Primera
$ws.Cells.Item(616, 1).Value = 1
$ws.Cells.Item(616, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(616, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(616, 4).Value = 44714
$ws.Cells.Item(616, 5).Value = 15
$ws.Cells.Item(616, 6).Value = 100112024
$ws.Cells.Item(616, 7).Value = "Choclo"
$ws.Cells.Item(616, 8).Value = "Lluteño"
$ws.Cells.Item(616, 9).Value = "Primera"
$ws.Cells.Item(616, 10).Value = 50
$ws.Cells.Item(616, 11).Value = 33000
$ws.Cells.Item(616, 12).Value = 35000
$ws.Cells.Item(616, 13).Value = 34000
$ws.Cells.Item(616, 14).Value = "$/saco 50 unidades"
$ws.Cells.Item(616, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(616, 16).Value = 680
$ws.Cells.Item(616, 17).Value = 50
$ws.Cells.Item(616, 18).Value = "Hortaliza"

# New row 617: Lluteño / Segunda
$ws.Cells.Item(617, 1).Value = 1
$ws.Cells.Item(617, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(617, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(617, 4).Value = 44714
$ws.Cells.Item(617, 5).Value = 15
$ws.Cells.Item(617, 6).Value = 100112024
$ws.Cells.Item(617, 7).Value = "Choclo"
$ws.Cells.Item(617, 8).Value = "Lluteño"
$ws.Cells.Item(617, 9).Value = "Segunda"
$ws.Cells.Item(617, 10).Value = 70
$ws.Cells.Item(617, 11).Value = 30000
$ws.Cells.Item(617, 12).Value = 32000
$ws.Cells.Item(617, 13).Value = 31000
$ws.Cells.Item(617, 14).Value = "$/saco 75 unidades"
$ws.Cells.Item(617, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(617, 16).Value = 413
$ws.Cells.Item(617, 17).Value = 75
$ws.Cells.Item(617, 18).Value = "Hortaliza"

# New row 618: Lluteño / Tercera
$ws.Cells.Item(618, 1).Value = 1
$ws.Cells.Item(618, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(618, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(618, 4).Value = 44714
$ws.Cells.Item(618, 5).Value = 15
$ws.Cells.Item(618, 6).Value = 100112024
$ws.Cells.Item(618, 7).Value = "Choclo"
$ws.Cells.Item(618, 8).Value = "Lluteño"
$ws.Cells.Item(618, 9).Value = "Tercera"
$ws.Cells.Item(618, 10).Value = 70
$ws.Cells.Item(618, 11).Value = 23000
$ws.Cells.Item(618, 12).Value = 25000
$ws.Cells.Item(618, 13).Value = 24000
$ws.Cells.Item(618, 14).Value = "$/saco 100 unidades"
$ws.Cells.Item(618, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(618, 16).Value = 240
$ws.Cells.Item(618, 17).Value = 100
$ws.Cells.Item(618, 18).Value = "Hortaliza"
